$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-10 17:48:38'
$ws.Range("I2").Value = '30.3 mm'
$ws.Range("E3").Value = '2026-02-10 17:48:41'
$ws.Range("I3").Value = '20.1 mm'
$ws.Range("O3").Value = '0.9 °C'
$ws.Range("E4").Value = '2026-02-10 17:48:43'
$ws.Range("O4").Value = '11.8 °C'
$ws.Range("E5").Value = '2026-02-10 17:48:46'
$ws.Range("I5").Value = '27.2 mm'
$ws.Range("E6").Value = '2026-02-10 17:48:48'
$ws.Range("J6").Value = '1004.4 hPa'
$ws.Range("E7").Value = '2026-02-10 17:48:51'
$ws.Range("J7").Value = '1004.8 hPa'
$ws.Range("O7").Value = '14.9 °C'
$ws.Range("E8").Value = '2026-02-10 17:48:53'
$ws.Range("J8").Value = '1004.7 hPa'
$ws.Range("O8").Value = '11.6 °C'
$ws.Range("E9").Value = '2026-02-10 17:48:56'
$ws.Range("E10").Value = '2026-02-10 17:48:59'
$ws.Range("O10").Value = '10.0 °C'
$ws.Range("E11").Value = '2026-02-10 17:49:01'
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '88%'
$ws.Range("O11").Value = '7.3 °C'
$ws.Range("E12").Value = '2026-02-10 17:49:03'
$ws.Range("E13").Value = '2026-02-10 17:49:05'
$ws.Range("J13").Value = '1007.0 hPa'
$ws.Range("O13").Value = '4.9 °C'
$ws.Range("E14").Value = '2026-02-10 17:49:08'
$ws.Range("O14").Value = '13.2 °C'
$ws.Range("E15").Value = '2026-02-10 17:49:10'
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = '92%'
$ws.Range("O15").Value = '8.9 °C'
$ws.Range("E16").Value = '2026-02-10 17:49:13'
$ws.Range("I16").Value = '21.6 mm'
$ws.Range("O16").Value = '1.1 °C'
$ws.Range("E17").Value = '2026-02-10 17:49:15'
$ws.Range("K17").Value = '6.4 MJ/m2'
$ws.Range("O17").Value = '4.5 °C'
$ws.Range("E18").Value = '2026-02-10 17:49:18'
$ws.Range("J18").Value = '1004.6 hPa'
$ws.Range("O18").Value = '10.1 °C'
$ws.Range("E19").Value = '2026-02-10 17:49:20'
$ws.Range("L19").Value = '25.6 km/h - 239º 17:29 TU'
$ws.Range("O19").Value = '6.2 °C'
$ws.Range("E20").Value = '2026-02-10 17:49:23'
$ws.Range("I20").Value = '3.5 mm'
$ws.Range("O20").Value = '0.2 °C'
$ws.Range("E21").Value = '2026-02-10 17:49:25'
$ws.Range("J21").Value = '1006.4 hPa'
$ws.Range("O21").Value = '7.0 °C'
$ws.Range("E22").Value = '2026-02-10 17:49:28'
$ws.Range("I22").Value = '7.5 mm'
$ws.Range("E23").Value = '2026-02-10 17:49:30'
$ws.Range("I23").Value = '21.5 mm'
$ws.Range("E24").Value = '2026-02-10 17:49:33'
$ws.Range("J24").Value = '1006.3 hPa'
$ws.Range("O24").Value = '11.1 °C'
$ws.Range("E25").Value = '2026-02-10 17:49:35'
$ws.Range("O25").Value = '1.5 °C'
$ws.Range("E26").Value = '2026-02-10 17:49:38'
$ws.Range("J26").Value = '1003.5 hPa'
$ws.Range("E27").Value = '2026-02-10 17:49:40'
$ws.Range("L27").Value = '41.8 km/h - 244º 17:01 TU'
$ws.Range("O27").Value = '1.1 °C'
$ws.Range("E28").Value = '2026-02-10 17:49:43'
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = '81%'
$ws.Range("J28").Value = '1004.5 hPa'
$ws.Range("O28").Value = '8.7 °C'
$ws.Range("E29").Value = '2026-02-10 17:49:45'
$ws.Range("E30").Value = '2026-02-10 17:49:48'
$ws.Range("J30").Value = '1004.5 hPa'
$ws.Range("E31").Value = '2026-02-10 17:49:50'
$ws.Range("J31").Value = '1003.7 hPa'
$ws.Range("E32").Value = '2026-02-10 17:49:53'
$ws.Range("O32").Value = '10.2 °C'
$ws.Range("E33").Value = '2026-02-10 17:49:55'
$ws.Range("J33").Value = '1006.8 hPa'
$ws.Range("O33").Value = '3.9 °C'
$ws.Range("E34").Value = '2026-02-10 17:49:58'
$ws.Range("K34").Value = '7.0 MJ/m2'
$ws.Range("E35").Value = '2026-02-10 17:50:00'
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = '70%'
$ws.Range("E36").Value = '2026-02-10 17:50:03'
$ws.Range("O36").Value = '9.9 °C'
$ws.Range("E37").Value = '2026-02-10 17:50:05'
$ws.Range("J37").Value = '1005.9 hPa'
$ws.Range("O37").Value = '6.4 °C'
$ws.Range("E38").Value = '2026-02-10 17:50:08'
$ws.Range("O38").Value = '10.5 °C'
$ws.Range("E39").Value = '2026-02-10 17:50:10'
$ws.Range("O39").Value = '1.7 °C'
$ws.Range("E40").Value = '2026-02-10 17:50:13'
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = '89%'
$ws.Range("I40").Value = '6.1 mm'
$ws.Range("J40").Value = '1007.3 hPa'
$ws.Range("E41").Value = '2026-02-10 17:50:15'
$ws.Range("J41").Value = '1004.9 hPa'
$ws.Range("E42").Value = '2026-02-10 17:50:18'
$ws.Range("O42").Value = '10.2 °C'
$ws.Range("E43").Value = '2026-02-10 17:50:20'
$ws.Range("O43").Value = '9.1 °C'
$ws.Range("E44").Value = '2026-02-10 17:50:22'
$ws.Range("I44").Value = '19.6 mm'
$ws.Range("E45").Value = '2026-02-10 17:50:25'
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = '93%'
$ws.Range("I45").Value = '25.9 mm'
$ws.Range("J45").Value = '1005.9 hPa'
$ws.Range("O45").Value = '6.3 °C'
$ws.Range("E46").Value = '2026-02-10 17:50:27'
$ws.Range("J46").Value = '1006.2 hPa'
$ws.Range("K46").Value = '8.9 MJ/m2'
$ws.Range("O46").Value = '13.9 °C'
